# Trade #44 closed at 2026-02-17 08:38:50 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate metrics after the new closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.52   # Current Capital
$summary.Range("B4").Value = -0.48     # Total P&L $
$summary.Range("B5").Value = -0.22     # Total P&L %
$summary.Range("B6").Value = 44        # Total Trades
$summary.Range("B8").Value = 20        # Losing Trades
$summary.Range("B9").Value = 34.09     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.52      # Capital
$status.Range("D4").Value = 44         # Trades
$status.Range("E4").Value = -0.48      # P&L $
$status.Range("F4").Value = -0.48      # P&L %
$status.Range("G4").Value = 34.09      # Win Rate %

# ---------------------------------------------------------------------------
# Helper to append the new trade #44 row to a trade log sheet
# ---------------------------------------------------------------------------
function Add-Trade44Row($ws) {
    $ws.Range("A45").Value = 44

    $ws.Range("B45").NumberFormat = "@"
    $ws.Range("B45").Value = "2026-02-17"

    $ws.Range("C45").NumberFormat = "@"
    $ws.Range("C45").Value = "08:38:44"

    $ws.Range("D45").Value = "MarketMaking"
    $ws.Range("E45").Value = "DOWN"
    $ws.Range("F45").Value = 0.28
    $ws.Range("G45").Value = 0.2
    $ws.Range("H45").Value = "CLOSED"
    $ws.Range("I45").Value = -28.5714
    $ws.Range("J45").Value = -0.08
    $ws.Range("K45").Value = 99.52
    $ws.Range("L45").Value = 0
    $ws.Range("M45").Value = 0
    $ws.Range("N45").Value = 0.6
    $ws.Range("O45").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P45").Value = "early_exit"
    $ws.Range("Q45").Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet: append trade #44
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade44Row $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet: append trade #44 (mirrors All Trades)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade44Row $marketMaking
